# Decrement column E (剩余) by 1 for every data row, except row 36 which
# remains unchanged, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($r, 5)   # column E
    $cell.Value2 = $cell.Value2 - 1
}
